# Updated cryptos list on Thu Mar 23 07:22:06 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.679.31'
$ws.Range('E2').Value = '  -2.17%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.752.62'
$ws.Range('E3').Value = '  -2.71%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.003'
$ws.Range('E4').Value = '  -0.12%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '323.66'
$ws.Range('E5').Value = '  -4.73%  '

$ws.Range('E6').Value = '  -0.15%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4314'
$ws.Range('E7').Value = '  -5.92%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3656'
$ws.Range('E8').Value = '  -3.89%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '45.19'
$ws.Range('E9').Value = '  -0.11%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07461'
$ws.Range('E10').Value = '  -1.64%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.119'
$ws.Range('E11').Value = '  -2.94%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.001'
$ws.Range('E12').Value = '  -0.18%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '21.58'
$ws.Range('E13').Value = '  -4.16%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.156'
$ws.Range('E14').Value = '  -3.09%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.242'
$ws.Range('E15').Value = '  -4.36%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.751.85'
$ws.Range('E16').Value = '  -2.93%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001066'
$ws.Range('E17').Value = '  -2.47%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '87.97'
$ws.Range('E18').Value = '  +7.92%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06198'
$ws.Range('E19').Value = '  -7.86%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.9989'
$ws.Range('E20').Value = '  -0.13%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.13'
$ws.Range('E21').Value = '  -1.88%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.150'
$ws.Range('E22').Value = '  -4.18%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.5291'
$ws.Range('E23').Value = '  -5.38%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '27.699.22'
$ws.Range('E24').Value = '  -2.09%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '11.62'
$ws.Range('E25').Value = '  -2.34%  '

$ws.Range('E26').Value = '  -4.10%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.57'
$ws.Range('E27').Value = '  -0.55%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '153.65'
$ws.Range('E28').Value = '  +0.44%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.362'
$ws.Range('E29').Value = '  -0.16%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.950.87'
$ws.Range('E30').Value = '  -2.97%  '

$ws.Range('B31').Value = 'BitcoinCash'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '127.43'
$ws.Range('E31').Value = '  -4.19%  '

$ws.Range('B32').Value = 'ImmutableX'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.216'
$ws.Range('E32').Value = '  -2.53%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.718'
$ws.Range('E33').Value = '  -2.40%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.09154'
$ws.Range('E34').Value = '  -3.64%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.620'
$ws.Range('E35').Value = '  -10.14%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '12.63'
$ws.Range('E36').Value = '  +4.31%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02308'
$ws.Range('E37').Value = '  -2.06%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.2159'
$ws.Range('E38').Value = '  -6.46%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.102'
$ws.Range('E39').Value = '  -3.33%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.6470'
$ws.Range('E40').Value = '  -2.35%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.06093'
$ws.Range('E41').Value = '  -4.00%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.195'
$ws.Range('E42').Value = '  -3.49%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.427'
$ws.Range('E43').Value = '  -3.94%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '7.950'
$ws.Range('E44').Value = '  -5.13%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.9991'
$ws.Range('E45').Value = '  -0.11%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '13.84'
$ws.Range('E46').Value = '  -3.07%  '

$ws.Range('B47').Value = 'PancakeSwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.745'
$ws.Range('E47').Value = '  -3.35%  '

$ws.Range('B48').Value = 'Decentraland'
$ws.Range('C48').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.5933'
$ws.Range('E48').Value = '  -3.17%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '125.79'
$ws.Range('E49').Value = '  -4.14%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.974'
$ws.Range('E50').Value = '  -3.18%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06898'
$ws.Range('E51').Value = '  -3.77%  '
